# Fruta / hortaliza, semanal
# Insert a new weekly record at row 175 (pushing the existing rows 175-191
# down to 176-192) for "Feria Lagunitas de Puerto Montt" / Piña / Caramelo.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 175..191 down to 176..192, leaving a blank row 175.
$ws.Rows("175").Insert()

# Populate the new row 175 with this week's record.
$ws.Range("A175").Value = 4
$ws.Range("B175").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C175").Value = "Los Lagos"
$ws.Range("D175").Value = 44578
$ws.Range("E175").Value = 10
$ws.Range("F175").Value = "Fruta"
$ws.Range("G175").Value = 100108
$ws.Range("H175").Value = "Tropicales y subtropicales"
$ws.Range("I175").Value = 100108005
$ws.Range("J175").Value = "Piña"
$ws.Range("K175").Value = "Caramelo"
$ws.Range("L175").Value = "Tercera"
$ws.Range("M175").Value = 120
$ws.Range("N175").Value = 19000
$ws.Range("O175").Value = 20000
$ws.Range("P175").Value = 19500
$ws.Range("Q175").Value = "$/caja 16 unidades"
$ws.Range("R175").Value = "Ecuador"
$ws.Range("S175").Value = 1219
$ws.Range("T175").Value = 16
